# Trade #121 closed at 2026-02-17 16:04:33 - unknown UNKNOWN +0.000%
#
# Applies:
#   - Summary sheet roll-up metrics (Current Capital, Total P&L $,
#     Total Trades, Winning Trades, Win Rate %)
#   - Strategy Status sheet roll-up row for "MarketMaking"
#   - Appends new trade row (#121) to "All Trades" and "MarketMaking" sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1198.85   # Current Capital
$summary.Range("B4").Value = -1.16    # Total P&L $
$summary.Range("B6").Value = 121      # Total Trades
$summary.Range("B7").Value = 45       # Winning Trades
$summary.Range("B9").Value = 37.19    # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 98.85     # Capital
$status.Range("D4").Value = 121       # Trades
$status.Range("E4").Value = -1.16     # P&L $
$status.Range("F4").Value = -1.15     # P&L %
$status.Range("G4").Value = 37.19     # Win Rate %

# ---------------------------------------------------------------------
# 3. New trade row (#121) appended as row 122 to both
#    "All Trades" and "MarketMaking" sheets (identical data)
# ---------------------------------------------------------------------
$tradeSheets = @("All Trades", "MarketMaking")

foreach ($sheetName in $tradeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 122

    $ws.Cells.Item($row, 1).Value = 121              # A: Trade #

    # Date column is stored as plain text in this workbook (not a real
    # Excel date). Force text formatting before assignment so Excel
    # doesn't coerce the string into a date serial, then restore the
    # cell's style to the sheet's default (unstyled) look afterwards.
    $ws.Range("B" + $row).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"     # B: Date
    $ws.Range("B" + $row).Style = $ws.Range("A1").Style

    $ws.Cells.Item($row, 3).Value = "16:04:26"       # C: Time

    $ws.Cells.Item($row, 4).Value = "MarketMaking"   # D: Strategy
    $ws.Cells.Item($row, 5).Value = "UP"             # E: Side
    $ws.Cells.Item($row, 6).Value = 0.97             # F: Entry Price
    $ws.Cells.Item($row, 7).Value = 0.98             # G: Exit Price
    $ws.Cells.Item($row, 8).Value = "CLOSED"         # H: Status
    $ws.Cells.Item($row, 9).Value = 1.0309           # I: P&L %
    $ws.Cells.Item($row, 10).Value = 0.01            # J: P&L $
    $ws.Cells.Item($row, 11).Value = 98.85           # K: Capital After
    $ws.Cells.Item($row, 12).Value = 0               # L: Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0               # M: Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6             # N: Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"  # O: Entry Reason
    $ws.Cells.Item($row, 16).Value = "early_exit"    # P: Exit Reason
    $ws.Cells.Item($row, 17).Value = 0.11            # Q: Duration (min)
}
